# Populate column C ("FolderID") with the Google Drive file IDs for every
# metro/indicator row (rows 2-25). Previously only a handful of rows had a
# FolderID value; now every data row gets one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$folderIds = @{
    2  = "1u9ADBRckg4thM9HLXFZQ2SOz2ZO1thWK"
    3  = "1oyFRf90eC-nDQfIqPV_BG_nGuTIo7qlH"
    4  = "1X1JRgveCQ507sgR7U8Vken-cGX-Sv30s"
    5  = "1fJCTjjOS69iWB_U0pCqxG66Atu6C_bue"
    6  = "1Xkv_V0Pbc3DlhYP_35zzygGsOauDZVWD"
    7  = "1C2H6FHwvkOgDy1Q2YavW0LSUXi15euRP"
    8  = "1iYpeqiSbydj4I1lzYaiuZ6dYk6RGtqbr"
    9  = "1yz6c0XbhGNkmStm3YGqrxuLb59uupm26"
    10 = "1tNBAS361kPQrpnaVOvnu60QjG_AyENho"
    11 = "1FuDmpAccQyOrDToV0Eo1771JuW_Ovp2-"
    12 = "1ZqprzzagswbxAqCkAs6R6YCM43GgE8wt"
    13 = "1N4VxaCnXz3vHdHXtreWIbJobaohat5xL"
    14 = "15v7rxQybr8NNnyDf4rVq4QQW4fTU_KYz"
    15 = "1fUmB9ai8JtCP2DYqjISGR6TcIvHHwQ0g"
    16 = "1InNn45yyXrinDV_KvlaPqF3Xl_GRvcPk"
    17 = "1rrZbHkaJ8ERIprZd0hw6gqcw_VWAQeC9"
    18 = "1S_WOLAvNYcLJm5iDGChFk9Abmtf25ag5"
    19 = "1O6_m9LCEsB61rjDwFmsfaqVRtw10thzu"
    20 = "1nFlj06ey-MMzZITUcJjw5EZhk1_YsFC_"
    21 = "1l97xlAgKuV76GIWVn5xarWbvP6nLkJfy"
    22 = "1fe0rO3gmXQHCipqSgA65q0YK_KXoEXz4"
    23 = "1C5grehiUyd_xOY4w0vj3sXuEncOLnB3D"
    24 = "1Kl8x3bLM6igZYjZCvLjFwfG9EFcGC5lr"
    25 = "1q9KeS9x0gfsqWhm0u-Bp2yRqySo3N19I"
}

# Write in the same order the source data arrived in (matches the shared
# string table ordering produced upstream): ascending by row, except the
# very last pair (24, 25) which were appended with 25 preceding 24.
$orderedRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,25,24)
foreach ($row in $orderedRows) {
    $ws.Cells.Item($row, 3).Value = $folderIds[$row]
}

# Move the active selection to where the author last left it.
$ws.Range("C24").Select()
